# commit: "error solve ifrs list"
# The financial figures in rows 2-9 of the company_list sheet were re-derived
# (values were wildly out of scale / duplicated across years) - overwrite
# each metric cell with the corrected figure, and blank out the handful of
# cells (per-row) that no longer carry a value after the fix. Rows 7-9 lose
# every metric column (D:AI) entirely, keeping only the code/market/name in
# columns A-C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---------------------------------------------------------------
$ws.Range("D2").Value = 1390
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = -12
$ws.Range("H2").Value = -12
$ws.Range("I2").Value = -12
$ws.Range("J2").ClearContents()
$ws.Range("K2").Value = 1387
$ws.Range("L2").Value = 619
$ws.Range("M2").Value = 767
$ws.Range("N2").Value = 767
$ws.Range("O2").ClearContents()
$ws.Range("P2").Value = 438
$ws.Range("Q2").Value = -28
$ws.Range("R2").Value = -51
$ws.Range("S2").Value = 20
$ws.Range("T2").Value = 19
$ws.Range("U2").Value = -47
$ws.Range("V2").Value = 466
$ws.Range("W2").Value = 0.2
$ws.Range("X2").Value = -0.9
$ws.Range("Y2").Value = -1.65
$ws.Range("Z2").Value = -0.92
$ws.Range("AA2").Value = 80.70999999999999
$ws.Range("AB2").Value = 44.31
$ws.Range("AC2").Value = -15
$ws.Range("AD2").Value = -61.52
$ws.Range("AE2").Value = 876
$ws.Range("AF2").Value = 1.04
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 87619710

# --- Row 3 ---------------------------------------------------------------
$ws.Range("D3").Value = 1426
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = 10
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = -5
$ws.Range("I3").Value = -5
$ws.Range("J3").ClearContents()
$ws.Range("K3").Value = 2174
$ws.Range("L3").Value = 405
$ws.Range("M3").Value = 1769
$ws.Range("N3").Value = 1769
$ws.Range("O3").ClearContents()
$ws.Range("P3").Value = 646
$ws.Range("Q3").Value = -92
$ws.Range("R3").Value = -206
$ws.Range("S3").Value = 835
$ws.Range("T3").Value = 14
$ws.Range("U3").Value = -105
$ws.Range("V3").Value = 297
$ws.Range("W3").Value = 0.73
$ws.Range("X3").Value = -0.36
$ws.Range("Y3").Value = -0.41
$ws.Range("Z3").Value = -0.29
$ws.Range("AA3").Value = 22.88
$ws.Range("AB3").Value = 152.65
$ws.Range("AC3").Value = -5
$ws.Range("AD3").Value = -415.91
$ws.Range("AE3").Value = 1370
$ws.Range("AF3").Value = 1.6
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 129124213

# --- Row 4 ---------------------------------------------------------------
$ws.Range("D4").Value = 1138
$ws.Range("E4").Value = -151
$ws.Range("F4").Value = -151
$ws.Range("G4").Value = -173
$ws.Range("H4").Value = -163
$ws.Range("I4").Value = -162
$ws.Range("J4").Value = -1
$ws.Range("K4").Value = 2233
$ws.Range("L4").Value = 490
$ws.Range("M4").Value = 1743
$ws.Range("N4").Value = 1738
$ws.Range("O4").Value = 5
$ws.Range("P4").Value = 646
$ws.Range("Q4").Value = -98
$ws.Range("R4").Value = -128
$ws.Range("S4").Value = -20
$ws.Range("T4").Value = 19
$ws.Range("U4").Value = -118
$ws.Range("V4").Value = 270
$ws.Range("W4").Value = -13.26
$ws.Range("X4").Value = -14.33
$ws.Range("Y4").Value = -9.23
$ws.Range("Z4").Value = -7.4
$ws.Range("AA4").Value = 28.11
$ws.Range("AB4").Value = 127.67
$ws.Range("AC4").Value = -125
$ws.Range("AD4").Value = -16.27
$ws.Range("AE4").Value = 1346
$ws.Range("AF4").Value = 1.52
$ws.Range("AG4").ClearContents()
$ws.Range("AH4").ClearContents()
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 129124213

# --- Row 5 ---------------------------------------------------------------
$ws.Range("D5").Value = 1096
$ws.Range("E5").Value = -216
$ws.Range("F5").Value = -216
$ws.Range("G5").Value = 1947
$ws.Range("H5").Value = 1434
$ws.Range("I5").Value = 1435
$ws.Range("J5").Value = -1
$ws.Range("K5").Value = 6240
$ws.Range("L5").Value = 3220
$ws.Range("M5").Value = 3019
$ws.Range("N5").Value = 3019
$ws.Range("O5").ClearContents()
$ws.Range("P5").Value = 646
$ws.Range("Q5").Value = 27
$ws.Range("R5").Value = -386
$ws.Range("S5").Value = 102
$ws.Range("T5").Value = 11
$ws.Range("U5").Value = 16
$ws.Range("V5").Value = 356
$ws.Range("W5").Value = -19.71
$ws.Range("X5").Value = 130.78
$ws.Range("Y5").Value = 60.32
$ws.Range("Z5").Value = 33.84
$ws.Range("AA5").Value = 106.66
$ws.Range("AB5").Value = 348.59
$ws.Range("AC5").Value = 1111
$ws.Range("AD5").Value = 1.04
$ws.Range("AE5").Value = 2339
$ws.Range("AF5").Value = 0.5
$ws.Range("AG5").ClearContents()
$ws.Range("AH5").ClearContents()
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 129124213

# --- Row 6 ---------------------------------------------------------------
$ws.Range("D6").Value = 1017
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 6
$ws.Range("G6").Value = -1246
$ws.Range("H6").Value = -953
$ws.Range("I6").Value = -953
$ws.Range("K6").Value = 3802
$ws.Range("L6").Value = 1761
$ws.Range("M6").Value = 2040
$ws.Range("N6").Value = 2040
$ws.Range("P6").Value = 646
$ws.Range("Q6").Value = 39
$ws.Range("R6").Value = -118
$ws.Range("S6").Value = 47
$ws.Range("T6").Value = 15
$ws.Range("U6").Value = 24
$ws.Range("V6").Value = 386
$ws.Range("W6").Value = 0.5600000000000001
$ws.Range("X6").Value = -93.72
$ws.Range("Y6").Value = -37.66
$ws.Range("Z6").Value = -18.98
$ws.Range("AA6").Value = 86.31
$ws.Range("AB6").Value = 194.81
$ws.Range("AC6").Value = -738
$ws.Range("AD6").Value = -1.38
$ws.Range("AE6").Value = 1580
$ws.Range("AF6").Value = 0.64
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 129124213

# --- Rows 7-9 --------------------------------------------------------------
# Every metric cell (D:AI) is removed, leaving only the A/B/C identifying
# columns (code, market, company name) in place.
$ws.Range("D7:AI7").ClearContents()
$ws.Range("D8:AI8").ClearContents()
$ws.Range("D9:AI9").ClearContents()
